$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: inpatient cost (S1) -----------------------------------------
# New raw build-up for the mean (B2) and its SE (C2)
$ws.Range("B2").Formula = "=896+556+251+158"
$ws.Range("C2").Formula = "=ROUND(B2*SQRT(((23376-16836)/3.92)^2-(((12358-6551)/3.92)^2+((9919-7141)/3.92)^2))/(20106-9454-8530),0)"

# D2 used to hold a placeholder note; it now cites the same source as the
# other rows, with the same wrap-text formatting already used by D3/D4.
$ws.Range("D3").Copy()
$ws.Range("D2").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D2").Value = "skinner2018healthcare"

# Remove the yellow highlight fill from A2/B2/C2 (now unused).
$ws.Range("A2").Interior.ColorIndex = -4142  # xlNone
$ws.Range("B2").Interior.ColorIndex = -4142
$ws.Range("C2").Interior.ColorIndex = -4142

# C2 becomes an integer (no-decimal) value like C3/C4.
$ws.Range("C2").NumberFormat = "0"

# --- Row 3: outpatient cost (P1/S2) -------------------------------------
$ws.Range("B3").Formula = "=1781+436+204+163"
$ws.Range("B3").NumberFormat = "#,##0"
$ws.Range("C3").Formula = "=ROUND(B3*SQRT(((23376-16836)/3.92)^2-(((12358-6551)/3.92)^2+((9919-7141)/3.92)^2))/(20106-9454-8530),0)"
$ws.Range("C3").Interior.ColorIndex = -4142  # xlNone, match the other cells' "no fill" xf
$ws.Range("C3").NumberFormat = "0"

# --- Row 4: outpatient cost (P2) ----------------------------------------
$ws.Range("B4").Formula = "=1781+436+204+163"
$ws.Range("C4").Formula = "=ROUND(B4*SQRT(((23376-16836)/3.92)^2-(((12358-6551)/3.92)^2+((9919-7141)/3.92)^2))/(20106-9454-8530),0)"
$ws.Range("C4").Interior.ColorIndex = -4142
$ws.Range("C4").NumberFormat = "0"

Write-Output "edits applied"
